$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) cells whose new values would otherwise be
# auto-converted to numbers by Excel, so they stay text like the rest of the column.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

# Row 2
$ws.Range('D2').Value = '71.040.20'
$ws.Range('E2').Value = '  +4.92%  '

# Row 3
$ws.Range('D3').Value = '2.619.45'
$ws.Range('E3').Value = '  +5.48%  '

# Row 4
$ws.Range('E4').Value = '  +0.03%  '

# Row 5
$ws.Range('D5').Value = '605.56'
$ws.Range('E5').Value = '  +3.10%  '

# Row 6
$ws.Range('D6').Value = '181.72'
$ws.Range('E6').Value = '  +3.82%  '

# Row 7
$ws.Range('E7').Value = '  -0.04%  '

# Row 8
$ws.Range('E8').Value = '  +1.97%  '

# Row 9
$ws.Range('D9').Value = '2.618.48'

# Row 10
$ws.Range('E10').Value = '  +14.50%  '

# Row 11
$ws.Range('E11').Value = '  +0.52%  '

# Row 12
$ws.Range('E12').Value = '  +4.52%  '

# Row 13
$ws.Range('D13').Value = '5.05'
$ws.Range('E13').Value = '  +2.02%  '

# Row 14
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.045.57'
$ws.Range('E14').Value = '  +3.76%  '

# Row 15
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '26.74'
$ws.Range('E15').Value = '  +5.86%  '

# Row 16
$ws.Range('E16').Value = '  +7.88%  '

# Row 17
$ws.Range('D17').Value = '71.051.60'
$ws.Range('E17').Value = '  +5.04%  '

# Row 18
$ws.Range('D18').Value = '2.600.42'
$ws.Range('E18').Value = '  +5.07%  '

# Row 19
$ws.Range('D19').Value = '382.83'
$ws.Range('E19').Value = '  +10.49%  '

# Row 20
$ws.Range('E20').Value = '  +6.63%  '

# Row 21
$ws.Range('D21').Value = '11.48'
$ws.Range('E21').Value = '  +6.34%  '

# Row 22
$ws.Range('D22').Value = '4.18'
$ws.Range('E22').Value = '  +2.18%  '

# Row 23
$ws.Range('D23').Value = '72.13'
$ws.Range('E23').Value = '  +2.02%  '

# Row 24
$ws.Range('D24').Value = '4.46'
$ws.Range('E24').Value = '  +6.55%  '

# Row 25
$ws.Range('E25').Value = '  -0.01%  '

# Row 26
$ws.Range('D26').Value = '1.87'
$ws.Range('E26').Value = '  +11.20%  '

# Row 27
$ws.Range('E27').Value = '  +10.29%  '

# Row 28
$ws.Range('D28').Value = '2.749.02'
$ws.Range('E28').Value = '  +5.28%  '

# Row 29
$ws.Range('E29').Value = '  +0.02%  '

# Row 30
$ws.Range('D30').Value = '0.0₃0951'
$ws.Range('E30').Value = '  +6.69%  '

# Row 31
$ws.Range('E31').Value = '  +6.54%  '

# Row 32
$ws.Range('E32').Value = '  +4.53%  '

# Row 33
$ws.Range('E33').Value = '  +7.04%  '

# Row 34
$ws.Range('E34').Value = '  +4.50%  '

# Row 35
$ws.Range('E35').Value = '  +0.09%  '

# Row 36
$ws.Range('D36').Value = '163.93'
$ws.Range('E36').Value = '  -0.30%  '

# Row 37
$ws.Range('E37').Value = '  -0.08%  '

# Row 38
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '1.94'
$ws.Range('E38').Value = '  +12.44%  '

# Row 39
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').Value = '19.22'
$ws.Range('E39').Value = '  +5.20%  '

# Row 40
$ws.Range('E40').Value = '  +1.71%  '

# Row 41
$ws.Range('E41').Value = '  +6.79%  '

# Row 42
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').Value = '5.07'
$ws.Range('E42').Value = '  +6.38%  '

# Row 43
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.05%  '

# Row 44
$ws.Range('D44').Value = '2.58'
$ws.Range('E44').Value = '  +8.99%  '

# Row 45
$ws.Range('E45').Value = '  +2.53%  '

# Row 46
$ws.Range('D46').Value = '40.10'
$ws.Range('E46').Value = '  +3.80%  '

# Row 47
$ws.Range('D47').Value = '153.94'
$ws.Range('E47').Value = '  +4.17%  '

# Row 48
$ws.Range('E48').Value = '  +4.11%  '

# Row 49
$ws.Range('D49').Value = '0.0₆0273'
$ws.Range('E49').Value = '  +7.66%  '

# Row 50
$ws.Range('D50').Value = '0.534'
$ws.Range('E50').Value = '  +4.71%  '

# Row 51
$ws.Range('E51').Value = '  +7.11%  '
